{"js": "// The document contains a single table of two-digit \u00d7 two-digit\n// multiplication equations (\"NN\u00d7NN=NNNN\"). This commit refreshes the\n// generated problem set, swapping each old equation's text for a newly\n// generated one. Every \"before\" equation string occurs exactly once in\n// the document, so a plain text search keyed on the exact old string\n// unambiguously identifies the run to update.\nconst replacements = [\n  [\"34\u00d773=2482\", \"66\u00d740=2640\"],\n  [\"17\u00d772=1224\", \"98\u00d729=2842\"],\n  [\"49\u00d712=588\", \"51\u00d770=3570\"],\n  [\"38\u00d738=1444\", \"82\u00d782=6724\"],\n  [\"63\u00d736=2268\", \"41\u00d759=2419\"],\n  [\"86\u00d783=7138\", \"85\u00d773=6205\"],\n  [\"50\u00d736=1800\", \"94\u00d734=3196\"],\n  [\"80\u00d781=6480\", \"38\u00d755=2090\"],\n  [\"19\u00d761=1159\", \"46\u00d736=1656\"],\n  [\"22\u00d740=880\", \"61\u00d766=4026\"],\n  [\"19\u00d727=513\", \"73\u00d736=2628\"],\n  [\"84\u00d721=1764\", \"94\u00d742=3948\"],\n  [\"33\u00d716=528\", \"98\u00d766=6468\"],\n  [\"44\u00d767=2948\", \"67\u00d750=3350\"],\n  [\"88\u00d775=6600\", \"27\u00d765=1755\"],\n  [\"75\u00d795=7125\", \"83\u00d778=6474\"],\n  [\"57\u00d769=3933\", \"74\u00d734=2516\"],\n  [\"32\u00d735=1120\", \"86\u00d791=7826\"],\n  [\"94\u00d714=1316\", \"23\u00d739=897\"],\n  [\"34\u00d779=2686\", \"87\u00d766=5742\"],\n  [\"79\u00d769=5451\", \"43\u00d717=731\"],\n  [\"19\u00d744=836\", \"13\u00d779=1027\"],\n  [\"80\u00d770=5600\", \"45\u00d721=945\"],\n  [\"44\u00d780=3520\", \"97\u00d769=6693\"],\n  [\"39\u00d784=3276\", \"16\u00d791=1456\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find equation text: ${oldText}`);\n  }\n\n  // Replace every occurrence (there is exactly one per the diff, but loop\n  // defensively in case a value were ever duplicated).\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single table of two-digit x two-digit\n# multiplication equations (\"NN x NN=NNNN\"). This commit refreshes the\n# generated problem set, swapping each old equation's text for a newly\n# generated one. Every \"before\" equation string occurs exactly once in\n# the document, so Find/Replace keyed on the exact old string unambiguously\n# identifies the text run to update.\n$pairs = @(\n    @(\"34\u00d773=2482\", \"66\u00d740=2640\"),\n    @(\"17\u00d772=1224\", \"98\u00d729=2842\"),\n    @(\"49\u00d712=588\", \"51\u00d770=3570\"),\n    @(\"38\u00d738=1444\", \"82\u00d782=6724\"),\n    @(\"63\u00d736=2268\", \"41\u00d759=2419\"),\n    @(\"86\u00d783=7138\", \"85\u00d773=6205\"),\n    @(\"50\u00d736=1800\", \"94\u00d734=3196\"),\n    @(\"80\u00d781=6480\", \"38\u00d755=2090\"),\n    @(\"19\u00d761=1159\", \"46\u00d736=1656\"),\n    @(\"22\u00d740=880\", \"61\u00d766=4026\"),\n    @(\"19\u00d727=513\", \"73\u00d736=2628\"),\n    @(\"84\u00d721=1764\", \"94\u00d742=3948\"),\n    @(\"33\u00d716=528\", \"98\u00d766=6468\"),\n    @(\"44\u00d767=2948\", \"67\u00d750=3350\"),\n    @(\"88\u00d775=6600\", \"27\u00d765=1755\"),\n    @(\"75\u00d795=7125\", \"83\u00d778=6474\"),\n    @(\"57\u00d769=3933\", \"74\u00d734=2516\"),\n    @(\"32\u00d735=1120\", \"86\u00d791=7826\"),\n    @(\"94\u00d714=1316\", \"23\u00d739=897\"),\n    @(\"34\u00d779=2686\", \"87\u00d766=5742\"),\n    @(\"79\u00d769=5451\", \"43\u00d717=731\"),\n    @(\"19\u00d744=836\", \"13\u00d779=1027\"),\n    @(\"80\u00d770=5600\", \"45\u00d721=945\"),\n    @(\"44\u00d780=3520\", \"97\u00d769=6693\"),\n    @(\"39\u00d784=3276\", \"16\u00d791=1456\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    # Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n    #   MatchSoundsLike, MatchAllWordForms, Forward, Wrap, Format,\n    #   ReplaceWith, Replace)\n    # Wrap=1 -> wdFindContinue, Replace=2 -> wdReplaceAll\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Could not find equation text: $oldText\"\n    }\n}\n"}
